# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) for the first data row
# (row 2) on both the "zh-cn" and "de-de" report sheets to reflect
# the newly generated handback timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 04:38:02"
$wsZhCn.Range("H2").Value = "2016-03-12 04:38:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 04:38:05"
$wsDeDe.Range("H2").Value = "2016-03-12 04:38:24"
